# [taekwon] 몬스터 spell_effect 변경
# 일반 몬스터 spell_effect 변경
# Adds a new "elite warrior" monster row (row 13) to the monster sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 20101
$ws.Range("B13").Value = "몬스터_전사_엘리트"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 150
$ws.Range("E13").Value = "warrior"
$ws.Range("F13").Value = 2000101

[void]$ws.Range("H17").Select()
